$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 9
$ws.Range("B10").Value = 49
$ws.Range("C10").Value = 51
$ws.Range("D10").Value = 53
$ws.Range("E10").Value = 54
$ws.Range("F10").Value = 52
$ws.Range("G10").Value = 50

# Copy column-A formatting (bold font, thin border, center/top aligned)
# from the row above down into the newly added A10 cell.
$ws.Range("A9").Copy()
$ws.Range("A10").PasteSpecial(-4122)
